$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    ,("OCCCOC(C)COc1ccccc1","NON-Mutagenic (Consensus score: 0.3)",0,"51730-94-0",0)
    ,("O=C(O)c1ccccc1C","NON-Mutagenic (Consensus score: 1)",0,"118-90-1",0)
    ,("Oc1c(cccc1C)C","NON-Mutagenic (Consensus score: 0.675)",0,"576-26-1",0)
    ,("c2cc1cc(ccc1c(c2)C)C","NON-Mutagenic (Consensus score: 1)",0,"575-43-9",0)
    ,("O=C(OCCC(C)C)C=Cc1ccccc1","NON-Mutagenic (Consensus score: 1)",0,"7779-65-9",0)
    ,("O=C(O)c1ccc(cc1)C","NON-Mutagenic (Consensus score: 1)",0,"99-94-5",0)
    ,("Oc2ccc(OC1OC(CO)C(O)C(O)C1(O))cc2","NON-Mutagenic (Consensus score: 0.4)",0,"497-76-7",0)
    ,("Oc1ccccc1C","NON-Mutagenic (Consensus score: 1)",0,"95-48-7",0)
    ,("c1ccc(cc1)C(C)C","NON-Mutagenic (Consensus score: 1)",0,"98-82-8",0)
    ,("OC(c1ccc(cc1)C)(C)C","NON-Mutagenic (Consensus score: 0.75)",0,"1197-01-9",0)
    ,("c1ccc2c(c1)ccc(c2C)C","NON-Mutagenic (Consensus score: 1)",0,"28804-88-8",0)
    ,("O=C(O)c1ccc(O)cc1(O)","NON-Mutagenic (Consensus score: 1)",0,"89-86-1",0)
    ,("O=C4c5c(O)cc(O)cc5(OC(c1ccc2OC(CO)C(Oc2(c1))c3ccc(O)c(OC)c3)C4(O))","NON-Mutagenic (Consensus score: 0.45)",0,"22888-70-6",0)
    ,("O=C(OCc1ccccc1)Cc2ccccc2","NON-Mutagenic (Consensus score: 1)",0,"102-16-9",0)
    ,("O=C(OCCOC)c1ccccc1(C(=O)OCCOC)","Mutagenic (Consensus score: 1)",1,"117-82-8",1)
    ,("O=C(c1cc(oc1C)C)C","NON-Mutagenic (Consensus score: 0.35)",1,"10599-70-9",0)
    ,("Oc1cc(cc(c1)C)C","NON-Mutagenic (Consensus score: 1)",0,"108-68-9",0)
    ,("Oc1ccc(cc1(O))C","NON-Mutagenic (Consensus score: 1)",0,"452-86-8",0)
    ,("O=Cc1ccccc1(C=O)","NON-Mutagenic (Consensus score: 0.525)",1,"643-79-8",0)
    ,("O=C(OCCC)c1ccc(O)cc1","NON-Mutagenic (Consensus score: 1)",0,"94-13-3",0)
    ,("Oc1ccc(cc1)C=Cc2cc(O)cc(O)c2","NON-Mutagenic (Consensus score: 0.325)",0,"501-36-0",0)
    ,("O=C(O)CCC(c1ccc(O)cc1)(c2ccc(O)cc2)C","NON-Mutagenic (Consensus score: 0.825)",0,"126-00-1",0)
    ,("OCCOc1ccccc1","NON-Mutagenic (Consensus score: 1)",0,"122-99-6",0)
    ,("O=C(O)c1cccc(c1)C","NON-Mutagenic (Consensus score: 1)",0,"99-04-7",0)
    ,("O=C(OCCCC)c1ccc(O)cc1","NON-Mutagenic (Consensus score: 0.825)",0,"94-26-8",0)
    ,("Oc1ccc(c(c1)C)C","NON-Mutagenic (Consensus score: 1)",0,"95-65-8",0)
    ,("c1ccc2c(c1)cccc2C","NON-Mutagenic (Consensus score: 1)",0,"1321-94-4",0)
    ,("Oc1cc(ccc1C)C","NON-Mutagenic (Consensus score: 1)",0,"95-87-4",0)
    ,("Oc1ccc(cc1C)C","NON-Mutagenic (Consensus score: 0.675)",1,"105-67-9",0)
    ,("O=C(OC(C)C)Cc1ccccc1","NON-Mutagenic (Consensus score: 1)",0,"4861-85-2",0)
    ,("O(c1cc(cc2OCOc12)CC=C)C","NON-Mutagenic (Consensus score: 1)",0,"607-91-0",0)
    ,("O=C(O)C=Cc1ccc(C=CC(=O)O)cc1","NON-Mutagenic (Consensus score: 1)",0,"16323-43-6",0)
    ,("OC(C)COc1ccccc1","NON-Mutagenic (Consensus score: 0.4)",0,"770-35-4",0)
    ,("O=C(OCCCCCC)c1ccccc1","NON-Mutagenic (Consensus score: 1)",0,"6789-88-4",0)
    ,("Oc1cccc(O)c1C","NON-Mutagenic (Consensus score: 1)",0,"608-25-3",0)
    ,("c1ccc2cc(ccc2(c1))C","NON-Mutagenic (Consensus score: 0.65)",0,"91-57-6",0)
    ,("O=C(OC(C(=O)c1ccccc1)c2ccccc2)C","NON-Mutagenic (Consensus score: 0.575)",0,"574-06-1",0)
    ,("C=Cc1cccc(c1)C","NON-Mutagenic (Consensus score: 0.9)",0,"25013-15-4",0)
    ,("c1ccc5c(c1)ccc3c5(cc2ccc(c4c2c3CC4)C)","Mutagenic (Consensus score: 1)",1,"56-49-5",1)
    ,("O=c1cc(-c2ccccc2)oc2ccc3ccccc3c12","Mutagenic (Consensus score: 0.3)",0,"6051-87-2",1)
    ,("O=c1c(-c2ccc(O)cc2)coc2cc(O)cc(O)c12","NON-Mutagenic (Consensus score: 0.675)",1,"446-72-0",0)
    ,("Cc1ccc2oc(=O)ccc2c1","NON-Mutagenic (Consensus score: 1)",0,"92-48-8",0)
)

$startRow = 559
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
}

# Column width adjustments
$ws.Columns("A").ColumnWidth = 49.14
$ws.Columns("B").ColumnWidth = 58.14

# Update view / selection to match the final state
$excel.ActiveWindow.ScrollRow = 582
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A602").Select()